$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 243; existing rows 243.. shift down by one
# (this grows the used range from A1:R351 to A1:R352).
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with the new weekly record.
$ws.Cells.Item(243, 1).Value = 8
$ws.Cells.Item(243, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(243, 3).Value = 'Coquimbo'
$ws.Cells.Item(243, 4).Value = 44992
$ws.Cells.Item(243, 5).Value = 4
$ws.Cells.Item(243, 6).Value = 100112031
$ws.Cells.Item(243, 7).Value = 'Poroto verde'
$ws.Cells.Item(243, 8).Value = 'Magnum'
$ws.Cells.Item(243, 9).Value = 'Primera'
$ws.Cells.Item(243, 10).Value = 340
$ws.Cells.Item(243, 11).Value = 21000
$ws.Cells.Item(243, 12).Value = 22000
$ws.Cells.Item(243, 13).Value = 21500
$ws.Cells.Item(243, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(243, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(243, 16).Value = 860
$ws.Cells.Item(243, 17).Value = 25
$ws.Cells.Item(243, 18).Value = 'Hortaliza'

# Match the date cell styling used by the rest of column D.
$ws.Cells.Item(243, 4).NumberFormat = $ws.Cells.Item(244, 4).NumberFormat()
